$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-5), columns A-G
$data = @(
    @(20, 26.9, 27, 25.4, 25.3, 99.09999999999999, 0.1),
    @(25, 26.9, 27, 25.4, 25.3, 99.09999999999999, 0.2),
    @(30, 26.9, 27, 25.4, 25.3, 99.2, 0.2),
    @(35, 26.9, 27, 25.4, 25.4, 99, 0.1)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}

# Set explicit column widths for columns A-G (values chosen so the
# resulting pixel-quantized OOXML width lands as close as possible to
# the target widths of 11.7109375 / 14.7109375 / 14.7109375 / 14.7109375 /
# 15.7109375 / 14.7109375 / 14.7109375)
$ws.Columns.Item(1).ColumnWidth = 10.85
$ws.Columns.Item(2).ColumnWidth = 13.85
$ws.Columns.Item(3).ColumnWidth = 13.85
$ws.Columns.Item(4).ColumnWidth = 13.85
$ws.Columns.Item(5).ColumnWidth = 14.85
$ws.Columns.Item(6).ColumnWidth = 13.85
$ws.Columns.Item(7).ColumnWidth = 13.85
